# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the data block (rows 64-65),
# pushing the previously existing rows 64-73 down to rows 66-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 64 (existing rows 64:73 shift down to 66:75)
$ws.Rows("64:65").Insert()

# New row 64: Terminal Hortofrutícola Agro Chillán - Espárragos, "Primera"
$row64 = New-Object 'object[,]' 1,18
$row64[0,0]  = 7
$row64[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row64[0,2]  = "Ñuble"
$row64[0,3]  = 45258
$row64[0,4]  = 16
$row64[0,5]  = 300000000
$row64[0,6]  = "Espárragos"
$row64[0,7]  = "Sin especificar"
$row64[0,8]  = "Primera"
$row64[0,9]  = 200
$row64[0,10] = 1500
$row64[0,11] = 1500
$row64[0,12] = 1500
$row64[0,13] = "`$/kilo"
$row64[0,14] = "Región de Ñuble"
$row64[0,15] = 1500
$row64[0,16] = 1
$row64[0,17] = "Hortaliza"
$ws.Range("A64:R64").Value = $row64

# New row 65: Terminal Hortofrutícola Agro Chillán - Espárragos, "Segunda"
$row65 = New-Object 'object[,]' 1,18
$row65[0,0]  = 7
$row65[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row65[0,2]  = "Ñuble"
$row65[0,3]  = 45258
$row65[0,4]  = 16
$row65[0,5]  = 300000000
$row65[0,6]  = "Espárragos"
$row65[0,7]  = "Sin especificar"
$row65[0,8]  = "Segunda"
$row65[0,9]  = 200
$row65[0,10] = 1300
$row65[0,11] = 1300
$row65[0,12] = 1300
$row65[0,13] = "`$/kilo"
$row65[0,14] = "Región de Ñuble"
$row65[0,15] = 1300
$row65[0,16] = 1
$row65[0,17] = "Hortaliza"
$ws.Range("A65:R65").Value = $row65

Write-Output "Inserted rows 64-65; dimension now A1:R75"
